$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (Initial Value), D (Final Value), E (Target Value)
# rows 2..26, reflecting the "select only adults" filter applied upstream.
$data = @{
    2  = @(1404.55, 2999.55, 2989.87)
    3  = @(194.79, 510.09, 1644.4285)
    4  = @(90.52, 118.44, 298.987)
    5  = @(35.89, 69.31999999999999, 448.4805000000001)
    6  = @(29.47, 53.1, 31)
    7  = @(122.22, 214.06, 300)
    8  = @(219.6, 1322.27, 868)
    9  = @(0.9, 2.16, 29.89870000000001)
    10 = @(11.01, 25.14, 298.987)
    11 = @(11.47, 16.5, 179.3922)
    12 = @(2623.22, 3053.74, 1)
    13 = @(2241.43, 3804.86, 3510)
    14 = @(13.55, 16.64, 6.8)
    15 = @(309.36, 532.41, 303)
    16 = @(0.7, 1.16, 0.9)
    17 = @(0.95, 1.25, 1)
    18 = @(0.52, 0.97, 1.1)
    19 = @(11.74, 15.53, 11.5)
    20 = @(4.06, 1.82, 2)
    21 = @(3.74, 56.88, 66.09999999999999)
    22 = @(41.23, 401.67, 560)
    23 = @(1.39, 1.85, 0.7)
    24 = @(481.27, 749.15, 322)
    25 = @(989.78, 1640.03, 649)
    26 = @(16.71, 14.77, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}
